$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2. Data reporter section: update organization / contact details.
# (Edited bottom-to-top, matching the order new shared strings were appended
# in the authoritative commit.)
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B9").Value = "(0312) 32 46 55"
$ws.Range("B8").Value = "yryskan.kalymbetova@gmail.com"
$ws.Range("B7").Value = "Kalymbetova Yryskan"
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz republic (Department of Household Statistics) under the UNICEF MICS global programme"

# Cells containing the new Cyrillic-authored text pick up a distinct font
# record (as Excel does when content is entered under a Cyrillic locale).
$ws.Range("B8").Font.Name = "Calibri"
$ws.Range("B9").Font.Name = "Calibri"
$ws.Range("B6").Font.Name = "Calibri"

# Update the active selection to reflect where the author left off editing.
$ws.Range("B8").Select()
